$wb = $excel.ActiveWorkbook

# ---- PIR sheet ----
$ws = $wb.Worksheets.Item("PIR")
$rows = @(
    @("2026-02-04", "14:14:48", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-04", "14:14:50", "14:00", "Bathroom", "Motion Detected", "Active"),
    @("2026-02-04", "14:14:54", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-04", "14:14:59", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-04", "14:15:04", "14:00", "Bathroom", "Motion Detected", "Active"),
    @("2026-02-04", "14:15:13", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-04", "14:15:18", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-04", "14:15:22", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-04", "14:15:28", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-04", "14:15:28", "14:00", "Bathroom", "Motion Detected", "Active"),
    @("2026-02-04", "14:15:35", "14:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-02-04", "14:15:38", "14:00", "Bathroom", "Motion Detected", "Active"),
    @("2026-02-04", "14:15:45", "14:00", "Bathroom", "No Motion", "Inactive")
)
$r = 165
foreach ($row in $rows) {
    $ws.Range("A" + $r).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# ---- Humidity sheet ----
$ws = $wb.Worksheets.Item("Humidity")
$rows = @(
    @("2026-02-04", "14:14:49", "14:00", "Bathroom", "76.7%", "Active"),
    @("2026-02-04", "14:14:51", "14:00", "Bathroom", "77.6%", "Active"),
    @("2026-02-04", "14:14:56", "14:00", "Bathroom", "76.7%", "Active"),
    @("2026-02-04", "14:15:01", "14:00", "Bathroom", "77.6%", "Active"),
    @("2026-02-04", "14:15:06", "14:00", "Bathroom", "76.6%", "Active"),
    @("2026-02-04", "14:15:12", "14:00", "Bathroom", "77.5%", "Active"),
    @("2026-02-04", "14:15:17", "14:00", "Bathroom", "76.6%", "Active"),
    @("2026-02-04", "14:15:27", "14:00", "Bathroom", "76.8%", "Active"),
    @("2026-02-04", "14:15:32", "14:00", "Bathroom", "77.7%", "Active"),
    @("2026-02-04", "14:15:37", "14:00", "Bathroom", "77.0%", "Active"),
    @("2026-02-04", "14:15:42", "14:00", "Bathroom", "78.0%", "Active"),
    @("2026-02-04", "14:15:47", "14:00", "Bathroom", "77.1%", "Active")
)
$r = 136
foreach ($row in $rows) {
    $ws.Range("A" + $r).NumberFormat = "@"
    $ws.Range("E" + $r).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

# ---- Temperature sheet ----
$ws = $wb.Worksheets.Item("Temperature")
$rows = @(
    @("2026-02-04", "14:14:49", "14:00", "Bathroom", "24.8C", "Active"),
    @("2026-02-04", "14:14:52", "14:00", "Bathroom", "24.8C", "Active"),
    @("2026-02-04", "14:14:57", "14:00", "Bathroom", "24.8C", "Active"),
    @("2026-02-04", "14:15:02", "14:00", "Bathroom", "24.8C", "Active"),
    @("2026-02-04", "14:15:07", "14:00", "Bathroom", "24.8C", "Active"),
    @("2026-02-04", "14:15:12", "14:00", "Bathroom", "24.8C", "Active"),
    @("2026-02-04", "14:15:17", "14:00", "Bathroom", "24.8C", "Active"),
    @("2026-02-04", "14:15:27", "14:00", "Bathroom", "24.8C", "Active"),
    @("2026-02-04", "14:15:32", "14:00", "Bathroom", "24.7C", "Active"),
    @("2026-02-04", "14:15:37", "14:00", "Bathroom", "24.7C", "Active"),
    @("2026-02-04", "14:15:42", "14:00", "Bathroom", "24.7C", "Active"),
    @("2026-02-04", "14:15:47", "14:00", "Bathroom", "24.7C", "Active")
)
$r = 136
foreach ($row in $rows) {
    $ws.Range("A" + $r).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

